# Update workbook per commit:
#  - Sheet name date changed from 20241126-090821 to 20241127-095436
#  - Column G (date values) on every data row changed from 45622 to 45623

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab
$ws.Name = "IClientBalance-20241127-095436-"

# Update the G column (date serial) for all data rows: 2..274
for ($r = 2; $r -le 274; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq 45622) {
        $cell.Value2 = 45623
    }
}
